# Scheduled-runner price/profit refresh for leve crafting sheets.
# For each affected row, H:N (current market price + computed NQ/HQ totals & profit)
# are refreshed with freshly pulled prices. Some rows gain/lose a cell when a recipe
# switches between NQ-only and HQ-available pricing (LevePriceNQ/HQ of 0 => blank profit col).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 40: Stuck in the Moment | Horn Glue
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

# Row 58: A Matter of Vital Importance | Mega-Potion of Vitality
$ws.Range("H58").Value = 4295.7144
$ws.Range("J58").Value = 6000
$ws.Range("L58").Value = 18000
$ws.Range("N58").Value = -18300

# Row 99: Rumor Has It | Commanding Craftsman's Tea
$ws.Range("H99").Value = 188.2
$ws.Range("I99").Value = 188.2
$ws.Range("K99").Value = 564.5999999999999
$ws.Range("M99").Value = 933.4000000000001

# Row 100: Asking for a Friend | Beetle Glue
$ws.Range("H100").Value = 1130.0588
$ws.Range("I100").Value = 614.13336
$ws.Range("J100").Value = 4999.5
$ws.Range("K100").Value = 614.13336
$ws.Range("L100").Value = 4999.5
$ws.Range("M100").Value = -73.13336000000004
$ws.Range("N100").Value = -6081.5

# Row 101: Edge of the Arcane | Cunning Craftsman's Tea
$ws.Range("H101").Value = 699.6667
$ws.Range("I101").Value = 899.5
$ws.Range("J101").Value = 300
$ws.Range("K101").Value = 2698.5
$ws.Range("L101").Value = 900
$ws.Range("M101").Value = -1076.5
$ws.Range("N101").Value = -4144

# Row 104: Pep-stepper | Infusion of Vitality
$ws.Range("H104").Value = 88.25
$ws.Range("I104").Value = 87
$ws.Range("J104").Value = 89.5
$ws.Range("K104").Value = 261
$ws.Range("L104").Value = 268.5
$ws.Range("M104").Value = 1486
$ws.Range("N104").Value = -3762.5

# Row 112: Making Ends Meet | Superior Spiritbond Potion
$ws.Range("H112").Value = 3033.3125
$ws.Range("J112").Value = 3393.6924
$ws.Range("L112").Value = 10181.0772
$ws.Range("N112").Value = -12397.0772

# Row 127: Liquid Competence | Competent Craftsman's Draught
$ws.Range("H127").Value = 843.5
$ws.Range("I127").Value = 843.5
$ws.Range("K127").Value = 2530.5
$ws.Range("M127").Value = 2429.5

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 6972.1562
$ws.Range("I138").Value = 8217
$ws.Range("J138").Value = 5873.7646
$ws.Range("K138").Value = 24651
$ws.Range("L138").Value = 17621.2938
$ws.Range("M138").Value = -19511
$ws.Range("N138").Value = -27901.2938


# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 4: Eyes Bigger than the Plate | Bronze Plate
$ws.Range("H4").Value = 1150
$ws.Range("I4").Value = 1150
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1150
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -1034

# Row 6: Don't Hit Me One More Time | Bronze Hoplon
$ws.Range("H6").Value = 59521
$ws.Range("J6").Value = 8311.538
$ws.Range("L6").Value = 8311.538
$ws.Range("N6").Value = -8657.538

# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 15544.338
$ws.Range("I32").Value = 6174.116
$ws.Range("J32").Value = 27394.912
$ws.Range("K32").Value = 6174.116
$ws.Range("L32").Value = 27394.912
$ws.Range("M32").Value = -5887.116
$ws.Range("N32").Value = -27968.912

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 3857.8958
$ws.Range("I132").Value = 1599.7354
$ws.Range("J132").Value = 9342
$ws.Range("K132").Value = 4799.206200000001
$ws.Range("L132").Value = 28026
$ws.Range("M132").Value = -2269.206200000001
$ws.Range("N132").Value = -33086


# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 4410.9644
$ws.Range("I105").Value = 3753.3125
$ws.Range("K105").Value = 3753.3125
$ws.Range("M105").Value = -2006.3125

# Row 107: The Gold Experience | Deepgold Nugget
$ws.Range("H107").Value = 2826.5264
$ws.Range("I107").Value = 1730.8462
$ws.Range("J107").Value = 5200.5
$ws.Range("K107").Value = 1730.8462
$ws.Range("L107").Value = 5200.5
$ws.Range("M107").Value = 189.1538
$ws.Range("N107").Value = -9040.5

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 3011.5908
$ws.Range("I134").Value = 2181.2778
$ws.Range("K134").Value = 6543.8334
$ws.Range("M134").Value = -4008.8334


# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 22: Driving Up the Wall | Elm Lumber
$ws.Range("H22").Value = 383.33334
$ws.Range("I22").Value = 150
$ws.Range("K22").Value = 150
$ws.Range("M22").Value = 200

# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 4713.1665
$ws.Range("I31").Value = 3967.182
$ws.Range("J31").Value = 5344.385
$ws.Range("K31").Value = 3967.182
$ws.Range("L31").Value = 5344.385
$ws.Range("M31").Value = -3672.182
$ws.Range("N31").Value = -5934.385

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 4713.1665
$ws.Range("I34").Value = 3967.182
$ws.Range("J34").Value = 5344.385
$ws.Range("K34").Value = 3967.182
$ws.Range("L34").Value = 5344.385
$ws.Range("M34").Value = -3765.182
$ws.Range("N34").Value = -5748.385

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 3799.1304
$ws.Range("I58").Value = 1748.6
$ws.Range("K58").Value = 1748.6
$ws.Range("M58").Value = -1545.6

# Row 105: Zelkova, My Love | Zelkova Lumber
$ws.Range("H105").Value = 3424.5454
$ws.Range("I105").Value = 3462.4285
$ws.Range("K105").Value = 3462.4285
$ws.Range("M105").Value = -1715.4285

# Row 109: Playing the Market | White Oak Necklace
$ws.Range("H109").Value = 19808.092
$ws.Range("J109").Value = 19808.092
$ws.Range("L109").Value = 19808.092
$ws.Range("N109").Value = -21888.092

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 3799.1304
$ws.Range("I136").Value = 1748.6
$ws.Range("K136").Value = 5245.799999999999
$ws.Range("M136").Value = -2695.799999999999


# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 4: In Hot Water | Boiled Egg
$ws.Range("H4").Value = 508.16666
$ws.Range("I4").Value = 496
$ws.Range("K4").Value = 1488
$ws.Range("M4").Value = -1376

# Row 36: Love's Crumpets Lost | Crumpet
$ws.Range("H36").Value = 2555.7144
$ws.Range("I36").Value = 481.66666
$ws.Range("K36").Value = 1444.99998
$ws.Range("M36").Value = -1275.99998

# Row 59: Comfort Me with Mushrooms | Buttons in a Blanket
$ws.Range("H59").Value = 8292.5
$ws.Range("I59").Value = 7477.5
$ws.Range("K59").Value = 22432.5
$ws.Range("M59").Value = -21892.5

# Row 81: It Goes Down Smoothly | Frozen Spirits
$ws.Range("H81").Value = 2586
$ws.Range("J81").Value = 2586
$ws.Range("L81").Value = 7758
$ws.Range("N81").Value = -10004

# Row 84: Quenching the Flame (L) | Frozen Spirits
$ws.Range("H84").Value = 2586
$ws.Range("J84").Value = 2586
$ws.Range("L84").Value = 23274
$ws.Range("N84").Value = -34506

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 1309
$ws.Range("I131").Value = 700
$ws.Range("J131").Value = 1461.25
$ws.Range("K131").Value = 2100
$ws.Range("L131").Value = 4383.75
$ws.Range("M131").Value = 2940
$ws.Range("N131").Value = -14463.75


# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 70: Sky Is the Limit | Mythrite Ingot
$ws.Range("H70").Value = 6421.1665
$ws.Range("I70").Value = 4514.75
$ws.Range("J70").Value = 7374.375
$ws.Range("K70").Value = 4514.75
$ws.Range("L70").Value = 7374.375
$ws.Range("M70").Value = -4244.75
$ws.Range("N70").Value = -7914.375

# Row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws.Range("H73").Value = 6421.1665
$ws.Range("I73").Value = 4514.75
$ws.Range("J73").Value = 7374.375
$ws.Range("K73").Value = 4514.75
$ws.Range("L73").Value = 7374.375
$ws.Range("M73").Value = -3578.75
$ws.Range("N73").Value = -9246.375

# Row 97: If I'd a Koppranickel for Every Time... | Koppranickel Ingot
$ws.Range("H97").Value = 1456
$ws.Range("I97").Value = 1540.1904
$ws.Range("K97").Value = 1540.1904
$ws.Range("M97").Value = -1044.1904

# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 2178.48
$ws.Range("I102").Value = 1486.7778
$ws.Range("K102").Value = 1486.7778
$ws.Range("M102").Value = 135.2221999999999


# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 16: Saddle Sore | Hard Leather
$ws.Range("H16").Value = 3624.8
$ws.Range("I16").Value = 3624.8
$ws.Range("K16").Value = 3624.8
$ws.Range("M16").Value = -3454.8

# Row 24: On Their Feet Again | Hard Leather Espadrilles
$ws.Range("H24").Value = 30625
$ws.Range("J24").Value = 30625
$ws.Range("L24").Value = 30625
$ws.Range("N24").Value = -31311

# Row 40: Best Served Toad | Toad Leather
$ws.Range("H40").Value = 4833
$ws.Range("I40").Value = 4249.5
$ws.Range("K40").Value = 4249.5
$ws.Range("M40").Value = -4113.5

# Row 61: Spelling Me Softly | Raptor Leather
$ws.Range("H61").Value = 3066.8572
$ws.Range("I61").Value = 2661.6667
$ws.Range("J61").Value = 5498
$ws.Range("K61").Value = 2661.6667
$ws.Range("L61").Value = 5498
$ws.Range("M61").Value = -2459.6667
$ws.Range("N61").Value = -5902

# Row 100: Tiger in the Sack | Tiger Leather
$ws.Range("H100").Value = 1572.8334
$ws.Range("I100").Value = 1487.4
$ws.Range("K100").Value = 1487.4
$ws.Range("M100").Value = -946.4000000000001

# Row 113: Peace in Rest | Atrociraptor Leather
$ws.Range("H113").Value = 3066.8572
$ws.Range("I113").Value = 2661.6667
$ws.Range("J113").Value = 5498
$ws.Range("K113").Value = 2661.6667
$ws.Range("L113").Value = 5498
$ws.Range("M113").Value = -491.6667000000002
$ws.Range("N113").Value = -9838

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 4548
$ws.Range("I132").Value = 3029.7778
$ws.Range("K132").Value = 9089.3334
$ws.Range("M132").Value = -6559.3334


# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 54: No Country for Cold Men | Woolen Tights
$ws.Range("H54").Value = 14166.667
$ws.Range("J54").Value = 14166.667
$ws.Range("L54").Value = 14166.667
$ws.Range("N54").Value = -15206.667

# Row 96: Skills on Display | Ruby Cotton Cloth
$ws.Range("H96").Value = 1832.6666
$ws.Range("J96").Value = 998
$ws.Range("L96").Value = 998
$ws.Range("N96").Value = -3744

# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 128374.875
$ws.Range("I126").Value = 500999.5
$ws.Range("J126").Value = 4166.6665
$ws.Range("K126").Value = 1502998.5
$ws.Range("L126").Value = 12499.9995
$ws.Range("M126").Value = -1500528.5
$ws.Range("N126").Value = -17439.9995

